# Resize the TestCaseTemplate columns per the June 10, 2024 docs update.
# Target "character" widths (as they appear in the workbook's column
# definitions) for columns A-F and I:
#   A: 23.29   B: 27.29   C: 41.71   D: 26.57   E: 40.43   F: 36.71   I: 25.29
#
# The host's ColumnWidth setter snaps to a pixel grid (1/6-character
# increments) before it is written back out, so we feed it the input that
# lands on the grid point closest to each target width rather than the
# raw target number itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 22.5
$ws.Columns.Item(2).ColumnWidth = 26.5
$ws.Columns.Item(3).ColumnWidth = 40.833333333333336
$ws.Columns.Item(4).ColumnWidth = 25.666666666666668
$ws.Columns.Item(5).ColumnWidth = 39.666666666666664
$ws.Columns.Item(6).ColumnWidth = 35.833333333333336
$ws.Columns.Item(9).ColumnWidth = 24.5
